# Reorders the weekly price records (rows 2-11) of the Chirimoya sheet.
# Only columns D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de
# comercializacion), S (Precio $/Kg) and T (Kg / unidad) are permuted
# across rows; all other columns (A, B, C, E-K, R) stay identical in
# every row, so they do not need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps each destination row to the source row whose values it should
# receive (both are 1-based worksheet row numbers).
$rowMap = @{
    2  = 4
    3  = 8
    4  = 5
    5  = 9
    6  = 3
    7  = 11
    8  = 7
    9  = 6
    10 = 10
    11 = 2
}

$cols = @("D", "L", "M", "N", "O", "P", "Q", "S", "T")

# Snapshot the original values of the affected columns for every row
# before writing anything back, since several rows both read from and
# are written to as part of the permutation.
$snapshot = @{}
foreach ($r in 2..11) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}
